# This deck has a single slide master/design ("Integral") that is applied to
# every slide, and its theme (ppt/theme/theme2.xml) drives the whole
# presentation's look. The edit being replicated here changes the deck's
# applied colour theme from the custom "Integral" palette over to the
# built-in "Office Theme" palette (the classic Office 2013+ colours), while
# everything else about the design (fonts, effects, layout) stays the same
# -- those two theme parts already shared an identical font/format scheme,
# only the 12 theme colours differ.
#
# We drive this through the slide's ThemeColorScheme, which is keyed by the
# standard MsoThemeColorSchemeIndex order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hyperlink, 12 followedHyperlink
# RGB() isn't available in this host, so each target colour's 0xBBGGRR OLE
# colour value is supplied directly (computed from the target hex triplets
# below). Because there is only one slide master shared by every slide, it
# is enough to apply this through slide 1 -- the whole presentation (every
# slide) uses that same master/theme.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# index : target colour name -> target hex -> OLE (BGR) value
$tcs.Item(1).RGB  = 0         # dk1      #000000
$tcs.Item(2).RGB  = 16777215  # lt1      #FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      #44546A
$tcs.Item(4).RGB  = 15132391  # lt2      #E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  #5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  #ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  #A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  #FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  #4472C4
$tcs.Item(10).RGB = 4697456   # accent6  #70AD47
$tcs.Item(11).RGB = 12673797  # hlink    #0563C1
$tcs.Item(12).RGB = 7491477   # folHlink #954F72
